# Wrap the two textbook titles in single quotation marks and prefix them
# with "the ", turning:
#   ... belongs to Textbook of Pathology by Harsh Mohan, Basic Pathology by ...
# into:
#   ... belongs to the 'Textbook of Pathology' by Harsh Mohan, the 'Basic Pathology' by ...
# (the quotes are the typographic U+2018/U+2019 marks used in the diff)

$d = $word.ActiveDocument
$lsquo = [char]0x2018
$rsquo = [char]0x2019

# --- "Textbook of Pathology" -> "the 'Textbook of Pathology'" ---
$r1 = $d.Content
$found1 = $r1.Find.Execute("Textbook of Pathology", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    # Insert the closing quote right after the matched phrase, then the
    # opening "the '" right before it. $r1.Start is unaffected by the
    # InsertAfter call, so it still points at the start of "Textbook".
    $r1.InsertAfter($rsquo)
    $startRng1 = $d.Range($r1.Start, $r1.Start)
    $startRng1.InsertBefore("the " + $lsquo)
}

# --- "Basic Pathology" -> "the 'Basic Pathology'" ---
$r2 = $d.Content
$found2 = $r2.Find.Execute("Basic Pathology", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $r2.InsertAfter($rsquo)
    $startRng2 = $d.Range($r2.Start, $r2.Start)
    $startRng2.InsertBefore("the " + $lsquo)
}
